$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "week 5" (SEMANA = 5) rows to append after the existing data block
# (row 113 is the last populated data row in the original sheet; rows
# 114-141 are currently empty in sheetData).
$newRows = @(
    @(2026, "ENERO", 5, "AMARILLO", "COLORES", "GOLDFINCH", 9189),
    @(2026, "ENERO", 5, "AMARILLO", "COLORES", "HIGH AND EXOTIC", 23455),
    @(2026, "ENERO", 5, "AMARILLO", "COLORES", "MOMENTUM", 10474),
    @(2026, "ENERO", 5, "AMARILLO", "COLORES", "SUNDAY MORNING", 655),
    @(2026, "ENERO", 5, "BIC. AMARILLO", "COLORES", "SUMMER LIGHT", 5700),
    @(2026, "ENERO", 5, "BICOLOR", "COLORES", "BLUSH", 2597),
    @(2026, "ENERO", 5, "BICOLOR", "COLORES", "DISCOVERY", 2250),
    @(2026, "ENERO", 5, "BLANCO", "COLORES", "HIGH AND PURE", 27237),
    @(2026, "ENERO", 5, "BLANCO", "COLORES", "SUGGAR DOLL", 6002),
    @(2026, "ENERO", 5, "BLANCO", "COLORES", "VANILLA ICE", 12230),
    @(2026, "ENERO", 5, "BLANCO", "COLORES", "VENDELA", 8320),
    @(2026, "ENERO", 5, "DURAZNO", "COLORES", "TIFANY", 5165),
    @(2026, "ENERO", 5, "HOT PINK", "COLORES", "COTTON CANDY", 7853),
    @(2026, "ENERO", 5, "HOT PINK", "COLORES", "JACARANDA", 16279),
    @(2026, "ENERO", 5, "HOT PINK", "COLORES", "PINK FLOYD", 14770),
    @(2026, "ENERO", 5, "LAVANDER", "COLORES", "DEEP PURPLE", 35110),
    @(2026, "ENERO", 5, "LAVANDER", "COLORES", "MOODY BLUES", 11882),
    @(2026, "ENERO", 5, "NARANJA", "COLORES", "ALIVE", 5370),
    @(2026, "ENERO", 5, "NARANJA", "COLORES", "BROMO", 4750),
    @(2026, "ENERO", 5, "NARANJA", "COLORES", "CLEMENTINA", 8020),
    @(2026, "ENERO", 5, "NARANJA", "COLORES", "NINA", 35775),
    @(2026, "ENERO", 5, "ROJO", "ROJO", "FREEDOM", 272061),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "ABSOLUT IN PINK", 5341),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "HIGH AND BONITA", 17550),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "LUCIANO", 2560),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "PECKOUBO", 1390),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "STARFISH", 6310),
    @(2026, "ENERO", 5, "ROSADO", "COLORES", "TABATHA", 12775)
)

$startRow = 114
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]
    $ws.Cells.Item($r, 3).Value = $rowValues[2]
    $ws.Cells.Item($r, 4).Value = $rowValues[3]
    $ws.Cells.Item($r, 5).Value = $rowValues[4]
    $ws.Cells.Item($r, 6).Value = $rowValues[5]
    $ws.Cells.Item($r, 7).Value = $rowValues[6]
}

# Update the active selection to match where the author left the cursor
# after pasting the new block (one row below the last new row).
[void]$ws.Range("A142").Select()
